$d = $word.ActiveDocument

# --- 1. Insert four new paragraphs right after the paragraph ending in
#        "...defending this." (currently followed by an existing blank
#        paragraph, then "PRACTICE REPEATEDLY:").
#        New paragraphs: blank, note #1, blank, note #2.

$findRng = $d.Content
$findRng.Find.Execute("defending this.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$total = $d.Paragraphs.Count
$anchorIndex = 0
for ($i = 1; $i -le $total; $i++) {
    $pr = $d.Paragraphs.Item($i).Range
    if ($pr.Start -le $findRng.Start -and $findRng.End -le $pr.End) {
        $anchorIndex = $i
    }
}

$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($anchorIndex + 1).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($anchorIndex + 2).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($anchorIndex + 3).Range.InsertParagraphAfter() | Out-Null

$d.Paragraphs.Item($anchorIndex + 2).Range.Text = "never say ""to be honest"" and don't use emotional language"
$d.Paragraphs.Item($anchorIndex + 4).Range.Text = "my contribution is methodological nature"

# --- 2. Flip the "Normal" style's overflow-punctuation setting from
#        false to true (w:overflowPunct).

$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.HangingPunctuation = $true
